# Apply the "tryb LEVEL i dalsze prace porzadkowe" edit to the Translation sheet.
#
# Summary of the change (per the supplied diff):
#  - Row 40, column F: "Autotest:"  ->  "Diagnostics:"
#  - Three new rows are appended to the Translation table:
#       Row 41: SingleUseId44 | Typography_00 | Center | LTR | <value>%
#       Row 42: SingleUseId45 | Typography_00 | Left   | LTR | MCU Load:
#       Row 43: SingleUseId46 | Typography_00 | Left   | LTR | Level:
#  - Row 44 is left as a blank spacer row (no data), matching the trailing
#    empty row present after the newly appended rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# --- Update the existing "Autotest:" row to read "Diagnostics:" ---
$ws.Cells.Item(40, 6).Value = "Diagnostics:"

# --- Append the new rows used by the "Level" diagnostics text ---

# Row 41
$ws.Cells.Item(41, 2).Value = "SingleUseId44"
$ws.Cells.Item(41, 3).Value = "Typography_00"
$ws.Cells.Item(41, 4).Value = "Center"
$ws.Cells.Item(41, 5).Value = "LTR"
$ws.Cells.Item(41, 6).Value = "<value>%"

# Row 42
$ws.Cells.Item(42, 2).Value = "SingleUseId45"
$ws.Cells.Item(42, 3).Value = "Typography_00"
$ws.Cells.Item(42, 4).Value = "Left"
$ws.Cells.Item(42, 5).Value = "LTR"
$ws.Cells.Item(42, 6).Value = "MCU Load: "

# Row 43
$ws.Cells.Item(43, 2).Value = "SingleUseId46"
$ws.Cells.Item(43, 3).Value = "Typography_00"
$ws.Cells.Item(43, 4).Value = "Left"
$ws.Cells.Item(43, 5).Value = "LTR"
$ws.Cells.Item(43, 6).Value = "Level:"

# Row 44 stays empty, acting as a trailing spacer row below the table.
# (Touching a no-op row property materializes the otherwise implicit
#  blank row without adding any visible formatting to it.)
$ws.Rows.Item(44).OutlineLevel = 0
